$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New daily COVID-19 Indonesia data rows (2020-04-03 .. 2020-04-06)
# Columns: A tanggal, B jumlah_periksa, C konfirmasi, D sembuh,
#          E meninggal, F negatif, G proses_periksa, H kasus_perawatan
$rows = @(
    @{ r = 47; tanggal = 43924; jumlah_periksa = 7986;  konfirmasi = 1986; sembuh = 134; meninggal = 181; negatif = 5715;  proses_periksa = 0; kasus_perawatan = 0 },
    @{ r = 48; tanggal = 43925; jumlah_periksa = 9712;  konfirmasi = 2092; sembuh = 150; meninggal = 191; negatif = 7620;  proses_periksa = 0; kasus_perawatan = 0 },
    @{ r = 49; tanggal = 43926; jumlah_periksa = 11242; konfirmasi = 2273; sembuh = 164; meninggal = 198; negatif = 8869;  proses_periksa = 0; kasus_perawatan = 0 },
    @{ r = 50; tanggal = 43927; jumlah_periksa = 13186; konfirmasi = 2491; sembuh = 192; meninggal = 209; negatif = 10695; proses_periksa = 0; kasus_perawatan = 0 }
)

# Column A carries the date number format (style index 2, "yyyy-mm-dd")
# already used by row 46 - clone it instead of assigning a NumberFormat
# string directly (which would mint a brand-new, differently-escaped
# number format / style entry).
$ws.Cells.Item(46, 1).Copy() | Out-Null

foreach ($row in $rows) {
    $r = $row.r
    $ws.Cells.Item($r, 1).PasteSpecial(-4122) | Out-Null
    $ws.Cells.Item($r, 1).Value = $row.tanggal
    $ws.Cells.Item($r, 2).Value = $row.jumlah_periksa
    $ws.Cells.Item($r, 3).Value = $row.konfirmasi
    $ws.Cells.Item($r, 4).Value = $row.sembuh
    $ws.Cells.Item($r, 5).Value = $row.meninggal
    $ws.Cells.Item($r, 6).Value = $row.negatif
    $ws.Cells.Item($r, 7).Value = $row.proses_periksa
    $ws.Cells.Item($r, 8).Value = $row.kasus_perawatan
}

$ws.Range("I50").Select()
